$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-120 (Generation 0-118) -> 7668
$ws.Range("C2:C120").Value = 7668

# Rows 121-140 (Generation 119-138) -> 7623
$ws.Range("C121:C140").Value = 7623

# Rows 141-252 (Generation 139-250) -> 7573
$ws.Range("C141:C252").Value = 7573
